# Documentation-Chess.docx: remove the "Iscrum" heading + icescrum link
# paragraphs, and refresh the cached footer page-number field (the page
# count shrinks from removing that content lower down, bringing the
# bookmarked closing paragraph onto the next page).

$d = $word.ActiveDocument

# --- 1. Locate the "Iscrum" heading paragraph and the paragraph that
#        holds the icescrum hyperlink, then delete the whole span
#        (heading + blank line + hyperlink line) while leaving the
#        trailing bookmark paragraph untouched.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Iscrum*") {
        $startPara = $p
        break
    }
}

$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*icescrum.cpnv.ch*") {
        $endPara = $p
        break
    }
}

if (($startPara -ne $null) -and ($endPara -ne $null)) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}

# --- 2. Update the cached PAGE field result in the footer from 2 to 3.
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)
